$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cadastro")

$ws1.Range("B2").Value = "john117"
$ws1.Range("B14").Value = "john117"
